$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.922976999999999
$ws.Range("H2").Value = 14.768931
$ws.Range("I2").Value = 0.2897120038548413
$ws.Range("J2").Value = 0.2897120038548412
$ws.Range("M2").Value = 6.875726333333334
$ws.Range("N2").Value = 20.627179
$ws.Range("O2").Value = 0.6245871044820662
$ws.Range("P2").Value = 0.6245871044820662
$ws.Range("Q2").Value = 33.84904259729433
$ws.Range("R2").Value = 304.641383375649
$ws.Range("S2").Value = 0.1809503816213925
$ws.Range("T2").Value = 0.1809503816213925

# Row 3
$ws.Range("G3").Value = 4.922976999999999
$ws.Range("H3").Value = 14.768931
$ws.Range("I3").Value = 0.2897120038548413
$ws.Range("J3").Value = 0.2897120038548412
$ws.Range("O3").Value = 0.06694469792011602
$ws.Range("P3").Value = 0.06694469792011602
$ws.Range("Q3").Value = 3.628019079004333
$ws.Range("R3").Value = 32.652171711039
$ws.Range("S3").Value = 0.01939468258189384
$ws.Range("T3").Value = 0.01939468258189384

# Row 4
$ws.Range("G4").Value = 4.922976999999999
$ws.Range("H4").Value = 14.768931
$ws.Range("I4").Value = 0.2897120038548413
$ws.Range("J4").Value = 0.2897120038548412
$ws.Range("O4").Value = 0.3084681975978177
$ws.Range("P4").Value = 0.3084681975978177
$ws.Range("Q4").Value = 16.717208993704
$ws.Range("R4").Value = 150.454880943336
$ws.Range("S4").Value = 0.08936693965155493
$ws.Range("T4").Value = 0.0893669396515549

# Row 5
$ws.Range("I5").Value = 0.5769489387710858
$ws.Range("J5").Value = 0.5769489387710858
$ws.Range("M5").Value = 6.875726333333334
$ws.Range("N5").Value = 20.627179
$ws.Range("O5").Value = 0.6245871044820662
$ws.Range("P5").Value = 0.6245871044820662
$ws.Range("Q5").Value = 67.40890589646133
$ws.Range("R5").Value = 606.6801530681521
$ws.Range("S5").Value = 0.3603548671010334
$ws.Range("T5").Value = 0.3603548671010334

# Row 6
$ws.Range("I6").Value = 0.5769489387710858
$ws.Range("J6").Value = 0.5769489387710858
$ws.Range("O6").Value = 0.06694469792011602
$ws.Range("P6").Value = 0.06694469792011602
$ws.Range("R6").Value = 65.025389236872
$ws.Range("S6").Value = 0.03862367242136185
$ws.Range("T6").Value = 0.03862367242136185

# Row 7
$ws.Range("I7").Value = 0.5769489387710858
$ws.Range("J7").Value = 0.5769489387710858
$ws.Range("O7").Value = 0.3084681975978177
$ws.Range("P7").Value = 0.3084681975978177
$ws.Range("R7").Value = 299.6243950481281
$ws.Range("S7").Value = 0.1779703992486905
$ws.Range("T7").Value = 0.1779703992486905

# Row 8
$ws.Range("I8").Value = 0.133339057374073
$ws.Range("J8").Value = 0.133339057374073
$ws.Range("M8").Value = 6.875726333333334
$ws.Range("N8").Value = 20.627179
$ws.Range("O8").Value = 0.6245871044820662
$ws.Range("P8").Value = 0.6245871044820662
$ws.Range("Q8").Value = 15.57891759017167
$ws.Range("R8").Value = 140.210258311545
$ws.Range("S8").Value = 0.08328185575964034
$ws.Range("T8").Value = 0.08328185575964034

# Row 9
$ws.Range("I9").Value = 0.133339057374073
$ws.Range("J9").Value = 0.133339057374073
$ws.Range("O9").Value = 0.06694469792011602
$ws.Range("P9").Value = 0.06694469792011602
$ws.Range("Q9").Value = 1.669784605721667
$ws.Range("S9").Value = 0.008926342916860334
$ws.Range("T9").Value = 0.008926342916860334

# Row 10
$ws.Range("I10").Value = 0.133339057374073
$ws.Range("J10").Value = 0.133339057374073
$ws.Range("O10").Value = 0.3084681975978177
$ws.Range("P10").Value = 0.3084681975978177
$ws.Range("Q10").Value = 7.694043945320002
$ws.Range("R10").Value = 69.24639550788001
$ws.Range("S10").Value = 0.0411308586975723
$ws.Range("T10").Value = 0.0411308586975723
